# Refresh the crypto price/volume table (GitHub Actions scheduled update).
# Updates Coin/Link/Price/Volume(1h) cells in columns B-E for rows 2-51;
# column A (rank index) and the header row are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "26.039.14"
$ws.Range("E2").Value = "  -1.16%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "1.642.40"
$ws.Range("E3").Value = "  -1.51%  "

# Row 4: TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.63%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.38"
$ws.Range("E5").Value = "  -0.93%  "

# Row 6: XRP
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5178"

# Row 7: USDC
$ws.Range("E7").Value = "  -0.55%  "

# Row 8: Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2613"
$ws.Range("E8").Value = "  -1.90%  "

# Row 9: Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06274"
$ws.Range("E9").Value = "  -2.13%  "

# Row 10: Solana
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.40"
$ws.Range("E10").Value = "  -1.93%  "

# Row 11: TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07744"
$ws.Range("E11").Value = "  -1.34%  "

# Row 12: WrappedEther
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.678.16"
$ws.Range("E12").Value = "  +0.79%  "

# Row 13: Polkadot
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.467"
$ws.Range("E13").Value = "  -2.24%  "

# Row 14: WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "1.867.18"
$ws.Range("E14").Value = "  -1.49%  "

# Row 15: Polygon
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5569"
$ws.Range("E15").Value = "  +0.47%  "

# Row 16: ShibaInu
$ws.Range("D16").Value = "0.0₅7981"
$ws.Range("E16").Value = "  -2.68%  "

# Row 17: Litecoin
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.65"

# Row 18: WrappedBTC
$ws.Range("D18").Value = "26.032.54"
$ws.Range("E18").Value = "  -1.26%  "

# Row 19: Dai
$ws.Range("E19").Value = "  -0.65%  "

# Row 20: Uniswap
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.614"
$ws.Range("E20").Value = "  -1.77%  "

# Row 21: BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "192.95"
$ws.Range("E21").Value = "  -0.36%  "

# Row 22: Avalanche
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.07"
$ws.Range("E22").Value = "  -2.32%  "

# Row 23: Chainlink
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.939"
$ws.Range("E23").Value = "  -1.91%  "

# Row 24: BinanceUSD
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.005"
$ws.Range("E24").Value = "  -0.62%  "

# Row 25: Monero
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.78"
$ws.Range("E25").Value = "  +0.30%  "

# Row 26: Stellar
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1200"
$ws.Range("E26").Value = "  -2.70%  "

# Row 27: Cosmos
$ws.Range("E27").Value = "  -1.13%  "

# Row 28: EthereumClassic
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.91"
$ws.Range("E28").Value = "  -1.24%  "

# Row 29: Toncoin
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.479"
$ws.Range("E29").Value = "  -1.34%  "

# Row 30: Hedera
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05640"
$ws.Range("E30").Value = "  -3.91%  "

# Row 31: PancakeSwap
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.263"
$ws.Range("E31").Value = "  -1.84%  "

# Row 32: InternetComputer(DFINITY)
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.450"
$ws.Range("E32").Value = "  -5.26%  "

# Row 33: Filecoin
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.356"
$ws.Range("E33").Value = "  +2.10%  "

# Row 34: LidoDAOToken
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.592"
$ws.Range("E34").Value = "  -0.98%  "

# Row 35: MXToken
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.789"
$ws.Range("E35").Value = "  -1.32%  "

# Row 36: HuobiToken
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.411"
$ws.Range("E36").Value = "  -0.37%  "

# Row 37: ARBITRUM
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9367"
$ws.Range("E37").Value = "  -3.61%  "

# Row 38: ImmutableX
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5656"
$ws.Range("E38").Value = "  -3.13%  "

# Row 39: FraxShare
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.951"
$ws.Range("E39").Value = "  +1.70%  "

# Row 40: VeChain
$ws.Range("E40").Value = "  -1.97%  "

# Row 41: Maker
$ws.Range("D41").Value = "1.051.46"
$ws.Range("E41").Value = "  -0.82%  "

# Row 42: PaxDollar
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.004"
$ws.Range("E42").Value = "  -0.64%  "

# Row 43: TrustWalletToken
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8386"
$ws.Range("E43").Value = "  -3.64%  "

# Row 44: Quant
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.28"
$ws.Range("E44").Value = "  -2.71%  "

# Row 45: RocketPoolETH
$ws.Range("D45").Value = "1.779.33"
$ws.Range("E45").Value = "  -1.49%  "

# Row 46: Aave
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.75"
$ws.Range("E46").Value = "  -2.00%  "

# Row 47: Frax
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.009"
$ws.Range("E47").Value = "  -0.67%  "

# Row 48: Cronos
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05311"
$ws.Range("E48").Value = "  +2.77%  "

# Row 49: BabyDogeCoin
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₈103"
$ws.Range("E49").Value = "  -2.93%  "

# Row 50: Mantle
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4326"
$ws.Range("E50").Value = "  -1.41%  "

# Row 51: EnergySwap
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.917"
$ws.Range("E51").Value = "  -1.27%  "
